$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove obsolete backlog items -------------------------------------
# "Create Test class for Unit Tests (researche)" (row 15) and
# "Create an alternative for closing the program via File -> Exit" (row 3)
# were dropped from the backlog. Delete the higher row first so the lower
# row number stays valid.
$ws.Rows("15").Delete()
$ws.Rows("3").Delete()

# --- Update remaining task texts & effort estimates ---------------------
# (row numbers below are the FINAL row numbers, after the two deletions)

$ws.Cells.Item(3, 4).Value = 10          # Create a class that will handle all Database related actions.

$ws.Cells.Item(4, 4).Value = 2.5         # Bind the references needed for the Excel database driver

$ws.Cells.Item(5, 4).Value = 3.5         # Implement a DataGridView for Displaying the Data pulled from the DB

$ws.Cells.Item(7, 4).Value = 2.5         # Implement a permanent saving of the DB path ...

$ws.Cells.Item(9, 4).Value = 2.5         # Create a new member user interface

$ws.Cells.Item(10, 4).Value = 5          # Handle all necessary CRUD operations with the DB ...

$ws.Cells.Item(13, 3).Value = "Implement the BackupManager, making a Backup before each change in the DB."

$ws.Cells.Item(12, 3).Value = "Build a member info (new payment) update section in the detailed information Window"

# --- Selection / cursor --------------------------------------------------
$ws.Range("D13").Select()
